$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H58").Value = 1226952.5
$ws.Range("J58").Value = 2971.4285
$ws.Range("L58").Value = 8914.2855
$ws.Range("N58").Value = -9214.2855

$ws.Range("H115").Value = 10540
$ws.Range("I115").Value = 14485.714
$ws.Range("J115").Value = 1333.3334
$ws.Range("K115").Value = 43457.142
$ws.Range("L115").Value = 4000.0002
$ws.Range("M115").Value = -41890.142
$ws.Range("N115").Value = -7134.0002

$ws.Range("H116").Value = 2149.8
$ws.Range("J116").Value = 1916.3334
$ws.Range("L116").Value = 1916.3334
$ws.Range("N116").Value = -8800.3334

$ws.Range("H137").Value = 2235.3635
$ws.Range("I137").Value = 2330
$ws.Range("J137").Value = 2121.8
$ws.Range("K137").Value = 6990
$ws.Range("L137").Value = 6365.400000000001
$ws.Range("M137").Value = -4440
$ws.Range("N137").Value = -11465.4

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 29444.182
$ws.Range("I32").Value = 13351.932
$ws.Range("J32").Value = 158182.19
$ws.Range("K32").Value = 13351.932
$ws.Range("L32").Value = 158182.19
$ws.Range("M32").Value = -13064.932
$ws.Range("N32").Value = -158756.19

$ws.Range("H44").Value = 12846.125
$ws.Range("J44").Value = 12824.143
$ws.Range("L44").Value = 12824.143
$ws.Range("N44").Value = -13800.143

$ws.Range("H61").Value = 1679.8379
$ws.Range("I61").Value = 1390.88
$ws.Range("J61").Value = 2281.8333
$ws.Range("K61").Value = 1390.88
$ws.Range("L61").Value = 2281.8333
$ws.Range("M61").Value = -1178.88
$ws.Range("N61").Value = -2705.8333

$ws.Range("H74").Value = 1519.3462
$ws.Range("I74").Value = 1435.95
$ws.Range("K74").Value = 1435.95
$ws.Range("M74").Value = -561.95

$ws.Range("H77").Value = 1519.3462
$ws.Range("I77").Value = 1435.95
$ws.Range("K77").Value = 7179.75
$ws.Range("M77").Value = -2811.75

$ws.Range("H132").Value = 10500.228
$ws.Range("I132").Value = 12168.127
$ws.Range("J132").Value = 2160.7273
$ws.Range("K132").Value = 36504.381
$ws.Range("L132").Value = 6482.1819
$ws.Range("M132").Value = -33974.381
$ws.Range("N132").Value = -11542.1819

$ws.Range("H136").Value = 1679.8379
$ws.Range("I136").Value = 1390.88
$ws.Range("J136").Value = 2281.8333
$ws.Range("K136").Value = 4172.64
$ws.Range("L136").Value = 6845.499899999999
$ws.Range("M136").Value = -1622.64
$ws.Range("N136").Value = -11945.4999

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H20").Value = 82453.84
$ws.Range("I20").Value = 96736.37
$ws.Range("J20").Value = 3900
$ws.Range("K20").Value = 96736.37
$ws.Range("L20").Value = 3900
$ws.Range("M20").Value = -96489.37
$ws.Range("N20").Value = -4394

$ws.Range("H99").Value = 1676.7727
$ws.Range("I99").Value = 1421.1875
$ws.Range("J99").Value = 2358.3333
$ws.Range("K99").Value = 1421.1875
$ws.Range("L99").Value = 2358.3333
$ws.Range("M99").Value = 76.8125
$ws.Range("N99").Value = -5354.3333

$ws.Range("H107").Value = 142926110
$ws.Range("I107").Value = 166743800
$ws.Range("K107").Value = 166743800
$ws.Range("M107").Value = -166741880

$ws.Range("H134").Value = 2947.025
$ws.Range("I134").Value = 2893.7463
$ws.Range("J134").Value = 3221.6155
$ws.Range("K134").Value = 8681.2389
$ws.Range("L134").Value = 9664.8465
$ws.Range("M134").Value = -6146.2389
$ws.Range("N134").Value = -14734.8465

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 28823.6
$ws.Range("I31").Value = 1438.2759
$ws.Range("K31").Value = 1438.2759
$ws.Range("M31").Value = -1143.2759

$ws.Range("H34").Value = 28823.6
$ws.Range("I34").Value = 1438.2759
$ws.Range("K34").Value = 1438.2759
$ws.Range("M34").Value = -1236.2759

$ws.Range("H51").Value = 7933.1113
$ws.Range("J51").Value = 7933.1113
$ws.Range("L51").Value = 7933.1113
$ws.Range("N51").Value = -9405.1113

$ws.Range("H61").Value = 7933.1113
$ws.Range("J61").Value = 7933.1113
$ws.Range("L61").Value = 7933.1113
$ws.Range("N61").Value = -8629.1113

$ws.Range("H94").Value = 1201.2916
$ws.Range("I94").Value = 1040.125
$ws.Range("J94").Value = 1281.875
$ws.Range("K94").Value = 1040.125
$ws.Range("L94").Value = 1281.875
$ws.Range("M94").Value = -589.125
$ws.Range("N94").Value = -2183.875

$ws.Range("H99").Value = 5914.0356
$ws.Range("I99").Value = 1956.5834
$ws.Range("J99").Value = 8882.125
$ws.Range("K99").Value = 1956.5834
$ws.Range("L99").Value = 8882.125
$ws.Range("M99").Value = -458.5834
$ws.Range("N99").Value = -11878.125

$ws.Range("H126").Value = 5914.0356
$ws.Range("I126").Value = 1956.5834
$ws.Range("J126").Value = 8882.125
$ws.Range("K126").Value = 5869.7502
$ws.Range("L126").Value = 26646.375
$ws.Range("M126").Value = -3399.7502
$ws.Range("N126").Value = -31586.375

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H17").Value = 2000
$ws.Range("I17").Value = 0
$ws.Range("J17").Value = 2000
$ws.Range("K17").Value = 0
$ws.Range("L17").Value = 6000
$ws.Range("M17").ClearContents()
$ws.Range("N17").Value = -6338

$ws.Range("H131").Value = 650146.5
$ws.Range("J131").Value = 745669.5
$ws.Range("L131").Value = 2237008.5
$ws.Range("N131").Value = -2247088.5

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H57").Value = 19800
$ws.Range("J57").Value = 19800
$ws.Range("L57").Value = 19800
$ws.Range("N57").Value = -21440

$ws.Range("H97").Value = 23256662
$ws.Range("I97").Value = 29412714
$ws.Range("J97").Value = 462.77777
$ws.Range("K97").Value = 29412714
$ws.Range("L97").Value = 462.77777
$ws.Range("M97").Value = -29412218
$ws.Range("N97").Value = -1454.77777

$ws.Range("H107").Value = 1010611.1
$ws.Range("I107").Value = 443
$ws.Range("J107").Value = 3367670
$ws.Range("K107").Value = 443
$ws.Range("L107").Value = 3367670
$ws.Range("M107").Value = 1477
$ws.Range("N107").Value = -3371510

$ws.Range("H113").Value = 1415.375
$ws.Range("I113").Value = 1127.75
$ws.Range("J113").Value = 1703
$ws.Range("K113").Value = 1127.75
$ws.Range("L113").Value = 1703
$ws.Range("M113").Value = 1042.25
$ws.Range("N113").Value = -6043

$ws.Range("H122").Value = 1976.6666
$ws.Range("I122").Value = 1660
$ws.Range("K122").Value = 4980
$ws.Range("M122").Value = -2530

$ws.Range("H126").Value = 4335.4443
$ws.Range("I126").Value = 4498.5
$ws.Range("J126").Value = 4009.3333
$ws.Range("K126").Value = 13495.5
$ws.Range("L126").Value = 12027.9999
$ws.Range("M126").Value = -11025.5
$ws.Range("N126").Value = -16967.9999

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 2621.4614
$ws.Range("I7").Value = 1615.8
$ws.Range("K7").Value = 1615.8
$ws.Range("M7").Value = -1503.8

$ws.Range("H46").Value = 3828.6843
$ws.Range("I46").Value = 4868.5713
$ws.Range("J46").Value = 3222.0833
$ws.Range("K46").Value = 4868.5713
$ws.Range("L46").Value = 3222.0833
$ws.Range("M46").Value = -4680.5713
$ws.Range("N46").Value = -3598.0833

$ws.Range("H93").Value = 3239.2727
$ws.Range("I93").Value = 3338.9412
$ws.Range("J93").Value = 2900.4
$ws.Range("K93").Value = 3338.9412
$ws.Range("L93").Value = 2900.4
$ws.Range("M93").Value = -2090.9412
$ws.Range("N93").Value = -5396.4

$ws.Range("H126").Value = 2621.4614
$ws.Range("I126").Value = 1615.8
$ws.Range("K126").Value = 4847.4
$ws.Range("M126").Value = -2377.4

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H81").Value = 182410.1
$ws.Range("I81").Value = 125516.5
$ws.Range("J81").Value = 334126.34
$ws.Range("K81").Value = 251033
$ws.Range("L81").Value = 668252.6800000001
$ws.Range("M81").Value = -249972
$ws.Range("N81").Value = -670374.6800000001

$ws.Range("H84").Value = 182410.1
$ws.Range("I84").Value = 125516.5
$ws.Range("J84").Value = 334126.34
$ws.Range("K84").Value = 1255165
$ws.Range("L84").Value = 3341263.4
$ws.Range("M84").Value = -1249861
$ws.Range("N84").Value = -3351871.4

$ws.Range("H127").Value = 27222.5
$ws.Range("J127").Value = 27222.5
$ws.Range("L127").Value = 27222.5
$ws.Range("N127").Value = -37142.5

$ws.Range("H136").Value = 2630.3635
$ws.Range("I136").Value = 3055.1667
$ws.Range("J136").Value = 2120.6
$ws.Range("K136").Value = 9165.500100000001
$ws.Range("L136").Value = 6361.799999999999
$ws.Range("M136").Value = -6615.500100000001
$ws.Range("N136").Value = -11461.8
